$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '52.032.27'
Set-TextValue 'E2' '  -0.48%  '
Set-TextValue 'D3' '2.959.96'
Set-TextValue 'E3' '  +2.42%  '
Set-TextValue 'E4' '  +0.11%  '
Set-TextValue 'D5' '353.64'
Set-TextValue 'E5' '  +0.03%  '
Set-TextValue 'D6' '112.41'
Set-TextValue 'E6' '  -0.58%  '
Set-TextValue 'E7' '  +1.20%  '
Set-TextValue 'E8' '  -0.01%  '
Set-TextValue 'D9' '0.634'
Set-TextValue 'E9' '  +1.67%  '
Set-TextValue 'D10' '39.72'
Set-TextValue 'E10' '  -2.47%  '
Set-TextValue 'D11' '0.0896'
Set-TextValue 'E11' '  +4.85%  '
Set-TextValue 'E12' '  +0.86%  '
Set-TextValue 'D13' '19.93'
Set-TextValue 'E13' '  -1.81%  '
Set-TextValue 'D14' '8.09'
Set-TextValue 'E14' '  +2.73%  '
Set-TextValue 'D15' '3.427.03'
Set-TextValue 'E15' '  +2.55%  '
Set-TextValue 'D16' '2.960.71'
Set-TextValue 'E16' '  +2.41%  '
Set-TextValue 'E17' '  +0.35%  '
Set-TextValue 'D18' '52.124.31'
Set-TextValue 'E18' '  -0.23%  '
Set-TextValue 'D19' '7.75'
Set-TextValue 'E19' '  +0.68%  '
Set-TextValue 'D20' '14.47'
Set-TextValue 'E20' '  +5.80%  '
Set-TextValue 'D21' '3.30'
Set-TextValue 'E21' '  -2.53%  '
Set-TextValue 'D22' '0.0₃0993'
Set-TextValue 'E22' '  +1.40%  '
Set-TextValue 'D23' '71.63'
Set-TextValue 'E23' '  +0.98%  '
Set-TextValue 'D24' '272.16'
Set-TextValue 'E24' '  +0.34%  '
Set-TextValue 'E25' '  +0.58%  '
Set-TextValue 'E26' '  +10.63%  '
Set-TextValue 'D27' '27.57'
Set-TextValue 'E27' '  +3.45%  '
Set-TextValue 'B28' 'Dai'
Set-TextValue 'C28' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D28' '0.999'
Set-TextValue 'E28' '  -0.11%  '
Set-TextValue 'B29' 'Filecoin'
Set-TextValue 'C29' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D29' '7.49'
Set-TextValue 'E29' '  +18.28%  '
Set-TextValue 'D30' '0.110'
Set-TextValue 'E30' '  +21.77%  '
Set-TextValue 'D31' '10.74'
Set-TextValue 'E31' '  +1.18%  '
Set-TextValue 'D32' '37.85'
Set-TextValue 'E32' '  -2.69%  '
Set-TextValue 'D33' '6.33'
Set-TextValue 'E33' '  +8.38%  '
Set-TextValue 'D34' '53.17'
Set-TextValue 'E34' '  +0.82%  '
Set-TextValue 'D35' '0.0452'
Set-TextValue 'E35' '  -0.60%  '
Set-TextValue 'B36' 'FirstDigitalUSD'
Set-TextValue 'C36' 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue 'D36' '0.998'
Set-TextValue 'E36' '  -0.09%  '
Set-TextValue 'B37' 'Toncoin'
Set-TextValue 'C37' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D37' '1.90'
Set-TextValue 'E37' '  -16.38%  '
Set-TextValue 'D38' '3.42'
Set-TextValue 'E38' '  +2.20%  '
Set-TextValue 'D39' '18.97'
Set-TextValue 'E39' '  -0.12%  '
Set-TextValue 'E40' '  +1.15%  '
Set-TextValue 'D41' '2.70'
Set-TextValue 'E41' '  +3.24%  '
Set-TextValue 'E42' '  +2.16%  '
Set-TextValue 'D43' '23.85'
Set-TextValue 'E43' '  +4.87%  '
Set-TextValue 'E44' '  -2.11%  '
Set-TextValue 'E45' '  +1.05%  '
Set-TextValue 'E46' '  +1.17%  '
Set-TextValue 'D47' '2.175.80'
Set-TextValue 'E47' '  -0.29%  '
Set-TextValue 'D48' '114.27'
Set-TextValue 'E48' '  -6.75%  '
Set-TextValue 'D49' '0.246'
Set-TextValue 'E49' '  +1.73%  '
Set-TextValue 'E50' '  +6.57%  '
Set-TextValue 'D51' '0.939'
Set-TextValue 'E51' '  -2.64%  '
